# "Created header for practice"
#
# The paragraph that used to read:
#     ** Extra practice **  Khan Academy offers interactive online units ...
# (a "First Paragraph" styled paragraph made of three runs: the
# "** Extra practice **" label, a single space, and the long body text)
# becomes a short bold "Extra practice" heading-like paragraph (still using
# the FirstParagraph style) followed by a new Body Text paragraph that
# holds the "Khan Academy ..." text.

$d = $word.ActiveDocument

# Locate the "** Extra practice ** " label together with the trailing
# space that separated it from the body text, and collapse the three runs
# down to the plain label text.
$rng = $d.Content
$found = $rng.Find.Execute("** Extra practice ** ", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'Extra practice' label to edit"
}

$rng.Text = "Extra practice"
$rng.Bold = 1

# Split the paragraph right after the label so "Extra practice" becomes its
# own paragraph and the Khan Academy sentence starts a new one.
$rng.Collapse(0)
$rng.InsertParagraphAfter()

# The newly created paragraph (containing the Khan Academy text) should use
# the Body Text style instead of inheriting First Paragraph.
$bodyRng = $d.Range($rng.End + 1, $rng.End + 1)
$bodyPara = $bodyRng.Paragraphs(1)
$bodyPara.Style = "BodyText"
